# Applies the "eosio_vs_eth" update:
#  - turns the header row into an AutoFilter table (A1:F23) + hidden _FilterDatabase name
#  - adds a new "Interface" comparison row (row 24) with 6 new shared strings
#  - updates sheet view (zoom + active cell selection)
#  - refreshes row heights to the values captured in the target workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. AutoFilter on the table header, plus the corresponding hidden defined name
#    (done before adding row 24 so the filter range stays A1:F23)
# ---------------------------------------------------------------------------
[void]$ws.Range("A1:F23").AutoFilter()
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$23")
$fd.Visible = $false

# ---------------------------------------------------------------------------
# 2. New row of data (row 24) - EOSIO vs ETH "Interface" comparison
# ---------------------------------------------------------------------------
$ws.Cells.Item(24, 1).Value = "Interface"
$ws.Cells.Item(24, 2).Value = @'
using distribute_action  = action_wrapper<"distribute"_n, &terraworlds::distribute>;
'@
$ws.Cells.Item(24, 3).Value = @'
interface InterfaceTerraworlds {
   function distribute(uint256 nextId) external returns (bool);
}
'@
$ws.Cells.Item(24, 4).Value = "Using contract's function from inside another contract."
$ws.Cells.Item(24, 5).Value = "Write the function inside the contract & make the function available to external contract(s) using action_wrapper"
$ws.Cells.Item(24, 6).Value = @'
Write the function inside the contract & make the function available via 2methods:
M-1: writing interface module inside another contract or 
M-2: creating a file with interface module inside.
'@

# carry over the same look (font/border/wrap) used by the rest of the table
$ws.Range("A23:F23").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(24).RowHeight = 153

# ---------------------------------------------------------------------------
# 3. Sheet view: zoom + active selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 135
$ws.Range("A22").Select()

# ---------------------------------------------------------------------------
# 4. Row heights refreshed to match current rendering
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 409.25
$ws.Rows.Item(3).RowHeight = 35.5
$ws.Rows.Item(4).RowHeight = 68
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 34
$ws.Rows.Item(7).RowHeight = 34
$ws.Rows.Item(8).RowHeight = 17
$ws.Rows.Item(9).RowHeight = 34
$ws.Rows.Item(10).RowHeight = 34
$ws.Rows.Item(11).RowHeight = 17
$ws.Rows.Item(12).RowHeight = 51
$ws.Rows.Item(13).RowHeight = 102
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 17
$ws.Rows.Item(16).RowHeight = 51
$ws.Rows.Item(17).RowHeight = 17
$ws.Rows.Item(18).RowHeight = 17
$ws.Rows.Item(19).RowHeight = 34
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(21).RowHeight = 119
$ws.Rows.Item(22).RowHeight = 119
$ws.Rows.Item(23).RowHeight = 85

Write-Output "edit complete"
